$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Attendance roll fixes -------------------------------------------------
# Two previously marked "excused absence" entries were corrected to
# "unexcused absence", and the 9/1 meeting row (row 9), which only had its
# date filled in, now has everyone's attendance recorded.

# Row 6: Younouss Thiam (H), 8/22 meeting: E -> U
$ws.Range("H6").Value = "U"

# Row 7: Brian Davis (I), 8/25 meeting: E -> U
$ws.Range("I7").Value = "U"

# Row 9: fill in attendance for the 9/1 meeting (previously blank).
$ws.Range("D9").Value = "A"
$ws.Range("E9").Value = "A"
$ws.Range("F9").Value = "A"
$ws.Range("G9").Value = "A"
$ws.Range("H9").Value = "U"
$ws.Range("I9").Value = "A"

# Leave the selection on H9, matching where the attendance was last edited.
$ws.Range("H9").Select()

$wb.Save()
